# Update the division-problem answers in the results table.
# Each populated row (1, 5, 9, 13, 17) has 5 answer cells that get new values.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1, Col 1: "39÷9=4, 3" -> "81÷7=11, 4"
$tbl.Cell(1, 1).Range.Text = "81÷7=11, 4"
# Row 1, Col 2: "46÷2=23, 0" -> "24÷4=6, 0"
$tbl.Cell(1, 2).Range.Text = "24÷4=6, 0"
# Row 1, Col 3: "31÷5=6, 1" -> "94÷9=10, 4"
$tbl.Cell(1, 3).Range.Text = "94÷9=10, 4"
# Row 1, Col 4: "88÷7=12, 4" -> "94÷4=23, 2"
$tbl.Cell(1, 4).Range.Text = "94÷4=23, 2"
# Row 1, Col 5: "25÷3=8, 1" -> "91÷6=15, 1"
$tbl.Cell(1, 5).Range.Text = "91÷6=15, 1"
# Row 5, Col 1: "71÷4=17, 3" -> "21÷2=10, 1"
$tbl.Cell(5, 1).Range.Text = "21÷2=10, 1"
# Row 5, Col 2: "26÷7=3, 5" -> "16÷6=2, 4"
$tbl.Cell(5, 2).Range.Text = "16÷6=2, 4"
# Row 5, Col 3: "16÷3=5, 1" -> "77÷3=25, 2"
$tbl.Cell(5, 3).Range.Text = "77÷3=25, 2"
# Row 5, Col 4: "21÷9=2, 3" -> "50÷4=12, 2"
$tbl.Cell(5, 4).Range.Text = "50÷4=12, 2"
# Row 5, Col 5: "35÷2=17, 1" -> "52÷9=5, 7"
$tbl.Cell(5, 5).Range.Text = "52÷9=5, 7"
# Row 9, Col 1: "52÷4=13, 0" -> "19÷9=2, 1"
$tbl.Cell(9, 1).Range.Text = "19÷9=2, 1"
# Row 9, Col 2: "64÷9=7, 1" -> "47÷5=9, 2"
$tbl.Cell(9, 2).Range.Text = "47÷5=9, 2"
# Row 9, Col 3: "44÷2=22, 0" -> "32÷3=10, 2"
$tbl.Cell(9, 3).Range.Text = "32÷3=10, 2"
# Row 9, Col 4: "74÷7=10, 4" -> "22÷3=7, 1"
$tbl.Cell(9, 4).Range.Text = "22÷3=7, 1"
# Row 9, Col 5: "64÷9=7, 1" -> "49÷2=24, 1"
$tbl.Cell(9, 5).Range.Text = "49÷2=24, 1"
# Row 13, Col 1: "16÷5=3, 1" -> "78÷7=11, 1"
$tbl.Cell(13, 1).Range.Text = "78÷7=11, 1"
# Row 13, Col 2: "34÷8=4, 2" -> "47÷9=5, 2"
$tbl.Cell(13, 2).Range.Text = "47÷9=5, 2"
# Row 13, Col 3: "32÷4=8, 0" -> "49÷7=7, 0"
$tbl.Cell(13, 3).Range.Text = "49÷7=7, 0"
# Row 13, Col 4: "84÷3=28, 0" -> "27÷2=13, 1"
$tbl.Cell(13, 4).Range.Text = "27÷2=13, 1"
# Row 13, Col 5: "36÷5=7, 1" -> "18÷3=6, 0"
$tbl.Cell(13, 5).Range.Text = "18÷3=6, 0"
# Row 17, Col 1: "22÷7=3, 1" -> "21÷8=2, 5"
$tbl.Cell(17, 1).Range.Text = "21÷8=2, 5"
# Row 17, Col 2: "63÷4=15, 3" -> "70÷3=23, 1"
$tbl.Cell(17, 2).Range.Text = "70÷3=23, 1"
# Row 17, Col 3: "64÷2=32, 0" -> "60÷9=6, 6"
$tbl.Cell(17, 3).Range.Text = "60÷9=6, 6"
# Row 17, Col 4: "92÷2=46, 0" -> "73÷8=9, 1"
$tbl.Cell(17, 4).Range.Text = "73÷8=9, 1"
# Row 17, Col 5: "94÷9=10, 4" -> "48÷3=16, 0"
$tbl.Cell(17, 5).Range.Text = "48÷3=16, 0"
